$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 135.36842
$ws.Range("I2").Value = 95.29412000000001
$ws.Range("K2").Value = 95.29412000000001
$ws.Range("M2").Value = 17.70587999999999
$ws.Range("H6").Value = 86071
$ws.Range("I6").Value = 86071
$ws.Range("K6").Value = 258213
$ws.Range("M6").Value = -258101
$ws.Range("H12").Value = 96.8
$ws.Range("I12").Value = 96.8
$ws.Range("K12").Value = 96.8
$ws.Range("M12").Value = 73.2
$ws.Range("H17").Value = 2420959
$ws.Range("J17").Value = 2469376.2
$ws.Range("L17").Value = 7408128.600000001
$ws.Range("N17").Value = -7408464.600000001
$ws.Range("H19").Value = 1940.05
$ws.Range("I19").Value = 961.0833
$ws.Range("K19").Value = 961.0833
$ws.Range("M19").Value = -786.0833
$ws.Range("H31").Value = 2309.3333
$ws.Range("I31").Value = 2309.3333
$ws.Range("K31").Value = 6927.999899999999
$ws.Range("M31").Value = -6697.999899999999
$ws.Range("H32").Value = 1991.4
$ws.Range("I32").Value = 1996.5
$ws.Range("J32").Value = 1988
$ws.Range("K32").Value = 1996.5
$ws.Range("L32").Value = 1988
$ws.Range("M32").Value = -1670.5
$ws.Range("N32").Value = -2640
$ws.Range("H40").Value = 18754244
$ws.Range("J40").Value = 37502710
$ws.Range("L40").Value = 37502710
$ws.Range("N40").Value = -37503060
$ws.Range("H41").Value = 786.375
$ws.Range("I41").Value = 239.2
$ws.Range("K41").Value = 239.2
$ws.Range("M41").Value = 200.8
$ws.Range("H51").Value = 9385.1875
$ws.Range("J51").Value = 9211.532999999999
$ws.Range("L51").Value = 9211.532999999999
$ws.Range("N51").Value = -10179.533
$ws.Range("H62").Value = 6673543.5
$ws.Range("I62").Value = 8340469
$ws.Range("J62").Value = 5840.2
$ws.Range("K62").Value = 8340469
$ws.Range("L62").Value = 5840.2
$ws.Range("M62").Value = -8339845
$ws.Range("N62").Value = -7088.2
$ws.Range("H64").Value = 2796.6785
$ws.Range("I64").Value = 2808.2
$ws.Range("J64").Value = 2700.6667
$ws.Range("K64").Value = 2808.2
$ws.Range("L64").Value = 2700.6667
$ws.Range("M64").Value = -2560.2
$ws.Range("N64").Value = -3196.6667
$ws.Range("H65").Value = 6673543.5
$ws.Range("I65").Value = 8340469
$ws.Range("J65").Value = 5840.2
$ws.Range("K65").Value = 41702345
$ws.Range("L65").Value = 29201
$ws.Range("M65").Value = -41699225
$ws.Range("N65").Value = -35441
$ws.Range("H67").Value = 2796.6785
$ws.Range("I67").Value = 2808.2
$ws.Range("J67").Value = 2700.6667
$ws.Range("K67").Value = 2808.2
$ws.Range("L67").Value = 2700.6667
$ws.Range("M67").Value = -1950.2
$ws.Range("N67").Value = -4416.6667
$ws.Range("H103").Value = 413.2
$ws.Range("I103").Value = 298.16666
$ws.Range("J103").Value = 489.8889
$ws.Range("K103").Value = 894.4999799999999
$ws.Range("L103").Value = 1469.6667
$ws.Range("M103").Value = -308.4999799999999
$ws.Range("N103").Value = -2641.6667
$ws.Range("H106").Value = 4213.2856
$ws.Range("I106").Value = 3665.6667
$ws.Range("J106").Value = 4624
$ws.Range("K106").Value = 3665.6667
$ws.Range("L106").Value = 4624
$ws.Range("M106").Value = -3034.6667
$ws.Range("N106").Value = -5886
$ws.Range("H113").Value = 5070.7856
$ws.Range("I113").Value = 5092.8184
$ws.Range("J113").Value = 4990
$ws.Range("K113").Value = 5092.8184
$ws.Range("L113").Value = 4990
$ws.Range("M113").Value = -1838.8184
$ws.Range("N113").Value = -11498
$ws.Range("H132").Value = 5208.8423
$ws.Range("I132").Value = 3377.0715
$ws.Range("K132").Value = 10131.2145
$ws.Range("M132").Value = -7601.2145
$ws.Range("H135").Value = 1606.8182
$ws.Range("I135").Value = 989.8
$ws.Range("K135").Value = 8908.199999999999
$ws.Range("M135").Value = -6373.199999999999
$ws.Range("H138").Value = 2928.0186
$ws.Range("J138").Value = 2794.6038
$ws.Range("L138").Value = 8383.811399999999
$ws.Range("N138").Value = -18663.8114
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4091.8
$ws.Range("I32").Value = 3498.1707
$ws.Range("K32").Value = 3498.1707
$ws.Range("M32").Value = -3211.1707
$ws.Range("H45").Value = 2835.3333
$ws.Range("I45").Value = 2842.5334
$ws.Range("K45").Value = 2842.5334
$ws.Range("M45").Value = -2465.5334
$ws.Range("H61").Value = 3919.611
$ws.Range("I61").Value = 2953.8215
$ws.Range("K61").Value = 2953.8215
$ws.Range("M61").Value = -2741.8215
$ws.Range("H74").Value = 69658.73
$ws.Range("I74").Value = 83152.25
$ws.Range("K74").Value = 83152.25
$ws.Range("M74").Value = -82278.25
$ws.Range("H77").Value = 69658.73
$ws.Range("I77").Value = 83152.25
$ws.Range("K77").Value = 415761.25
$ws.Range("M77").Value = -411393.25
$ws.Range("H88").Value = 2902.5
$ws.Range("I88").Value = 2970.75
$ws.Range("J88").Value = 2879.75
$ws.Range("K88").Value = 2970.75
$ws.Range("L88").Value = 2879.75
$ws.Range("M88").Value = -2564.75
$ws.Range("N88").Value = -3691.75
$ws.Range("H91").Value = 2902.5
$ws.Range("I91").Value = 2970.75
$ws.Range("J91").Value = 2879.75
$ws.Range("K91").Value = 2970.75
$ws.Range("L91").Value = 2879.75
$ws.Range("M91").Value = -1566.75
$ws.Range("N91").Value = -5687.75
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()
$ws.Range("H110").Value = 3606.0454
$ws.Range("I110").Value = 3491
$ws.Range("K110").Value = 3491
$ws.Range("M110").Value = -1446
$ws.Range("H132").Value = 3346.8462
$ws.Range("I132").Value = 3207.9167
$ws.Range("K132").Value = 9623.750100000001
$ws.Range("M132").Value = -7093.750100000001
$ws.Range("H136").Value = 3919.611
$ws.Range("I136").Value = 2953.8215
$ws.Range("K136").Value = 8861.4645
$ws.Range("M136").Value = -6311.4645
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H5").Value = 0
$ws.Range("J5").Value = 0
$ws.Range("L5").Value = 0
$ws.Range("N5").ClearContents()
$ws.Range("H12").Value = 447
$ws.Range("I12").Value = 198
$ws.Range("J12").Value = 571.5
$ws.Range("K12").Value = 198
$ws.Range("L12").Value = 571.5
$ws.Range("M12").Value = -30
$ws.Range("N12").Value = -907.5
$ws.Range("H20").Value = 3092.2
$ws.Range("I20").Value = 3092.2
$ws.Range("K20").Value = 3092.2
$ws.Range("M20").Value = -2845.2
$ws.Range("H22").Value = 638.5
$ws.Range("I22").Value = 638.5
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 638.5
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -465.5
$ws.Range("N22").ClearContents()
$ws.Range("H81").Value = 31980.834
$ws.Range("J81").Value = 28377
$ws.Range("L81").Value = 28377
$ws.Range("N81").Value = -30499
$ws.Range("H82").Value = 125044296
$ws.Range("I82").Value = 333339330
$ws.Range("K82").Value = 333339330
$ws.Range("M82").Value = -333338947
$ws.Range("H84").Value = 31980.834
$ws.Range("J84").Value = 28377
$ws.Range("L84").Value = 85131
$ws.Range("N84").Value = -95739
$ws.Range("H85").Value = 125044296
$ws.Range("I85").Value = 333339330
$ws.Range("K85").Value = 333339330
$ws.Range("M85").Value = -333338004
$ws.Range("H86").Value = 8779
$ws.Range("I86").Value = 8075.6665
$ws.Range("K86").Value = 8075.6665
$ws.Range("M86").Value = -6952.6665
$ws.Range("H89").Value = 8779
$ws.Range("I89").Value = 8075.6665
$ws.Range("K89").Value = 40378.3325
$ws.Range("M89").Value = -34762.3325
$ws.Range("H94").Value = 936.2222
$ws.Range("I94").Value = 761.5833
$ws.Range("J94").Value = 2333.3333
$ws.Range("K94").Value = 761.5833
$ws.Range("L94").Value = 2333.3333
$ws.Range("M94").Value = -310.5833
$ws.Range("N94").Value = -3235.3333
$ws.Range("H99").Value = 4251.933
$ws.Range("I99").Value = 3588.4614
$ws.Range("K99").Value = 3588.4614
$ws.Range("M99").Value = -2090.4614
$ws.Range("H105").Value = 1814.8889
$ws.Range("I105").Value = 1754.5518
$ws.Range("J105").Value = 2064.8572
$ws.Range("K105").Value = 1754.5518
$ws.Range("L105").Value = 2064.8572
$ws.Range("M105").Value = -7.551799999999957
$ws.Range("N105").Value = -5558.8572
$ws.Range("H107").Value = 752
$ws.Range("I107").Value = 742.1667
$ws.Range("K107").Value = 742.1667
$ws.Range("M107").Value = 1177.8333
$ws.Range("H112").Value = 50000
$ws.Range("J112").Value = 50000
$ws.Range("L112").Value = 50000
$ws.Range("N112").Value = -52954
$ws.Range("H134").Value = 3752.6382
$ws.Range("I134").Value = 4714.6
$ws.Range("K134").Value = 14143.8
$ws.Range("M134").Value = -11608.8
$ws.Range("H138").Value = 65366.4
$ws.Range("J138").Value = 65366.4
$ws.Range("L138").Value = 65366.4
$ws.Range("N138").Value = -75646.39999999999
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 466.44446
$ws.Range("I22").Value = 466.44446
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 466.44446
$ws.Range("L22").Value = 0
$ws.Range("M22").Value = -116.44446
$ws.Range("N22").ClearContents()
$ws.Range("H31").Value = 17120414
$ws.Range("I31").Value = 401570.7
$ws.Range("K31").Value = 401570.7
$ws.Range("M31").Value = -401275.7
$ws.Range("H34").Value = 17120414
$ws.Range("I34").Value = 401570.7
$ws.Range("K34").Value = 401570.7
$ws.Range("M34").Value = -401368.7
$ws.Range("H58").Value = 4630
$ws.Range("I58").Value = 3000
$ws.Range("J58").Value = 5037.5
$ws.Range("K58").Value = 3000
$ws.Range("L58").Value = 5037.5
$ws.Range("M58").Value = -2797
$ws.Range("N58").Value = -5443.5
$ws.Range("H86").Value = 34006470
$ws.Range("J86").Value = 9264.666999999999
$ws.Range("L86").Value = 9264.666999999999
$ws.Range("N86").Value = -11510.667
$ws.Range("H89").Value = 34006470
$ws.Range("J89").Value = 9264.666999999999
$ws.Range("L89").Value = 46323.335
$ws.Range("N89").Value = -57555.335
$ws.Range("H122").Value = 1963.3334
$ws.Range("I122").Value = 1956.2
$ws.Range("J122").Value = 1999
$ws.Range("K122").Value = 5868.6
$ws.Range("L122").Value = 5997
$ws.Range("M122").Value = -3418.6
$ws.Range("N122").Value = -10897
$ws.Range("H132").Value = 9179.049999999999
$ws.Range("I132").Value = 9281.272000000001
$ws.Range("J132").Value = 9054.111000000001
$ws.Range("K132").Value = 27843.816
$ws.Range("L132").Value = 27162.333
$ws.Range("M132").Value = -25313.816
$ws.Range("N132").Value = -32222.333
$ws.Range("H134").Value = 4762.8335
$ws.Range("I134").Value = 3715.4
$ws.Range("K134").Value = 11146.2
$ws.Range("M134").Value = -8611.200000000001
$ws.Range("H136").Value = 4630
$ws.Range("I136").Value = 3000
$ws.Range("J136").Value = 5037.5
$ws.Range("K136").Value = 9000
$ws.Range("L136").Value = 15112.5
$ws.Range("M136").Value = -6450
$ws.Range("N136").Value = -20212.5
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3148644.8
$ws.Range("I4").Value = 3577921
$ws.Range("J4").Value = 787625
$ws.Range("K4").Value = 10733763
$ws.Range("L4").Value = 2362875
$ws.Range("M4").Value = -10733651
$ws.Range("N4").Value = -2363099
$ws.Range("H14").Value = 1389712.8
$ws.Range("I14").Value = 1389712.8
$ws.Range("K14").Value = 4169138.4
$ws.Range("M14").Value = -4168965.4
$ws.Range("H25").Value = 190
$ws.Range("I25").Value = 190
$ws.Range("K25").Value = 570
$ws.Range("M25").Value = -401
$ws.Range("H30").Value = 190
$ws.Range("I30").Value = 190
$ws.Range("K30").Value = 570
$ws.Range("M30").Value = -468
$ws.Range("H34").Value = 319.6
$ws.Range("I34").Value = 93.55556
$ws.Range("J34").Value = 2354
$ws.Range("K34").Value = 280.66668
$ws.Range("L34").Value = 7062
$ws.Range("M34").Value = -196.66668
$ws.Range("N34").Value = -7230
$ws.Range("H36").Value = 1225
$ws.Range("I36").Value = 1225
$ws.Range("K36").Value = 3675
$ws.Range("M36").Value = -3506
$ws.Range("H68").Value = 2978127.5
$ws.Range("J68").Value = 2180.261
$ws.Range("L68").Value = 6540.782999999999
$ws.Range("N68").Value = -8162.782999999999
$ws.Range("H71").Value = 2978127.5
$ws.Range("J71").Value = 2180.261
$ws.Range("L71").Value = 19622.349
$ws.Range("N71").Value = -27734.349
$ws.Range("H129").Value = 36459710
$ws.Range("I129").Value = 48612290
$ws.Range("J129").Value = 1995.5
$ws.Range("K129").Value = 145836870
$ws.Range("L129").Value = 5986.5
$ws.Range("M129").Value = -145831870
$ws.Range("N129").Value = -15986.5
$ws.Range("H132").Value = 8208.111000000001
$ws.Range("I132").Value = 13641
$ws.Range("J132").Value = 2775.2222
$ws.Range("K132").Value = 122769
$ws.Range("L132").Value = 24976.9998
$ws.Range("M132").Value = -120239
$ws.Range("N132").Value = -30036.9998
$ws.Range("H133").Value = 4248.3335
$ws.Range("I133").Value = 2996.6667
$ws.Range("K133").Value = 8990.000100000001
$ws.Range("M133").Value = -3930.000100000001
$ws.Range("H139").Value = 2510.1333
$ws.Range("I139").Value = 1826.4584
$ws.Range("K139").Value = 5479.3752
$ws.Range("M139").Value = -339.3752000000004
$ws.Range("H140").Value = 11728
$ws.Range("I140").Value = 17751.75
$ws.Range("K140").Value = 53255.25
$ws.Range("M140").Value = -48075.25
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H11").Value = 867892.2
$ws.Range("I11").Value = 1555980.6
$ws.Range("J11").Value = 54696.727
$ws.Range("K11").Value = 1555980.6
$ws.Range("L11").Value = 54696.727
$ws.Range("M11").Value = -1555841.6
$ws.Range("N11").Value = -54974.727
$ws.Range("H14").Value = 9430002
$ws.Range("I14").Value = 14668846
$ws.Range("J14").Value = 82.2
$ws.Range("K14").Value = 14668846
$ws.Range("L14").Value = 82.2
$ws.Range("M14").Value = -14668678
$ws.Range("N14").Value = -418.2
$ws.Range("H18").Value = 99998.336
$ws.Range("I18").Value = 99999
$ws.Range("J18").Value = 99998
$ws.Range("K18").Value = 99999
$ws.Range("L18").Value = 99998
$ws.Range("M18").Value = -99706
$ws.Range("N18").Value = -100584
$ws.Range("H43").Value = 15588.286
$ws.Range("I43").Value = 10931.637
$ws.Range("J43").Value = 32662.666
$ws.Range("K43").Value = 10931.637
$ws.Range("L43").Value = 32662.666
$ws.Range("M43").Value = -10780.637
$ws.Range("N43").Value = -32964.666
$ws.Range("H46").Value = 19990
$ws.Range("J46").Value = 19990
$ws.Range("L46").Value = 19990
$ws.Range("N46").Value = -20302
$ws.Range("H70").Value = 166676830
$ws.Range("J70").Value = 250011250
$ws.Range("L70").Value = 250011250
$ws.Range("N70").Value = -250011790
$ws.Range("H73").Value = 166676830
$ws.Range("J73").Value = 250011250
$ws.Range("L73").Value = 250011250
$ws.Range("N73").Value = -250013122
$ws.Range("H102").Value = 51214.383
$ws.Range("J102").Value = 147378.28
$ws.Range("L102").Value = 147378.28
$ws.Range("N102").Value = -150622.28
$ws.Range("H122").Value = 3198.9
$ws.Range("I122").Value = 3171.2222
$ws.Range("J122").Value = 3448
$ws.Range("K122").Value = 9513.6666
$ws.Range("L122").Value = 10344
$ws.Range("M122").Value = -7063.6666
$ws.Range("N122").Value = -15244
$ws.Range("H126").Value = 9746.75
$ws.Range("I126").Value = 14600
$ws.Range("K126").Value = 43800
$ws.Range("M126").Value = -41330
$ws.Range("H132").Value = 92945.37
$ws.Range("I132").Value = 144315.28
$ws.Range("J132").Value = 3048
$ws.Range("K132").Value = 432945.84
$ws.Range("L132").Value = 9144
$ws.Range("M132").Value = -430415.84
$ws.Range("N132").Value = -14204
$ws.Range("H136").Value = 23262.234
$ws.Range("J136").Value = 23262.234
$ws.Range("L136").Value = 69786.702
$ws.Range("N136").Value = -74886.702
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3669.7932
$ws.Range("I7").Value = 3816.111
$ws.Range("J7").Value = 3603.95
$ws.Range("K7").Value = 3816.111
$ws.Range("L7").Value = 3603.95
$ws.Range("M7").Value = -3704.111
$ws.Range("N7").Value = -3827.95
$ws.Range("H13").Value = 2995
$ws.Range("I13").Value = 2995
$ws.Range("K13").Value = 2995
$ws.Range("M13").Value = -2855
$ws.Range("H20").Value = 7978.3335
$ws.Range("I20").Value = 7920
$ws.Range("K20").Value = 7920
$ws.Range("M20").Value = -7694
$ws.Range("H38").Value = 0
$ws.Range("I38").Value = 0
$ws.Range("J38").Value = 0
$ws.Range("K38").Value = 0
$ws.Range("L38").Value = 0
$ws.Range("M38").ClearContents()
$ws.Range("N38").ClearContents()
$ws.Range("H40").Value = 3077.4285
$ws.Range("I40").Value = 3077.4285
$ws.Range("K40").Value = 3077.4285
$ws.Range("M40").Value = -2941.4285
$ws.Range("H43").Value = 17208.092
$ws.Range("J43").Value = 17270.117
$ws.Range("L43").Value = 17270.117
$ws.Range("N43").Value = -17656.117
$ws.Range("H46").Value = 3200.926
$ws.Range("I46").Value = 2665.0527
$ws.Range("J46").Value = 4473.625
$ws.Range("K46").Value = 2665.0527
$ws.Range("L46").Value = 4473.625
$ws.Range("M46").Value = -2477.0527
$ws.Range("N46").Value = -4849.625
$ws.Range("H93").Value = 2739.5
$ws.Range("J93").Value = 4332.6665
$ws.Range("L93").Value = 4332.6665
$ws.Range("N93").Value = -6828.6665
$ws.Range("H100").Value = 16598.777
$ws.Range("I100").Value = 18298.625
$ws.Range("J100").Value = 3000
$ws.Range("K100").Value = 18298.625
$ws.Range("L100").Value = 3000
$ws.Range("M100").Value = -17757.625
$ws.Range("N100").Value = -4082
$ws.Range("H122").Value = 2918.5
$ws.Range("I122").Value = 2883.24
$ws.Range("J122").Value = 3800
$ws.Range("K122").Value = 8649.719999999999
$ws.Range("L122").Value = 11400
$ws.Range("M122").Value = -6199.719999999999
$ws.Range("N122").Value = -16300
$ws.Range("H126").Value = 3669.7932
$ws.Range("I126").Value = 3816.111
$ws.Range("J126").Value = 3603.95
$ws.Range("K126").Value = 11448.333
$ws.Range("L126").Value = 10811.85
$ws.Range("M126").Value = -8978.332999999999
$ws.Range("N126").Value = -15751.85
$ws.Range("H132").Value = 4190.25
$ws.Range("I132").Value = 4070.9412
$ws.Range("J132").Value = 4866.3335
$ws.Range("K132").Value = 12212.8236
$ws.Range("L132").Value = 14599.0005
$ws.Range("M132").Value = -9682.8236
$ws.Range("N132").Value = -19659.0005
$ws.Range("H136").Value = 3817.7297
$ws.Range("I136").Value = 3860.8333
$ws.Range("J136").Value = 3738.1538
$ws.Range("K136").Value = 11582.4999
$ws.Range("L136").Value = 11214.4614
$ws.Range("M136").Value = -9032.499899999999
$ws.Range("N136").Value = -16314.4614
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H6").Value = 8334.333000000001
$ws.Range("J6").Value = 8334.333000000001
$ws.Range("L6").Value = 8334.333000000001
$ws.Range("N6").Value = -8564.333000000001
$ws.Range("H9").Value = 3999.5
$ws.Range("I9").Value = 3000
$ws.Range("J9").Value = 4999
$ws.Range("K9").Value = 3000
$ws.Range("L9").Value = 4999
$ws.Range("M9").Value = -2860
$ws.Range("N9").Value = -5279
$ws.Range("H25").Value = 2000
$ws.Range("J25").Value = 2000
$ws.Range("L25").Value = 2000
$ws.Range("N25").Value = -2586
$ws.Range("H26").Value = 9500
$ws.Range("J26").Value = 9500
$ws.Range("L26").Value = 9500
$ws.Range("N26").Value = -10086
$ws.Range("H48").Value = 49999
$ws.Range("J48").Value = 49999
$ws.Range("L48").Value = 49999
$ws.Range("N48").Value = -51137
$ws.Range("H51").Value = 13998.6
$ws.Range("I51").Value = 11248.25
$ws.Range("J51").Value = 25000
$ws.Range("K51").Value = 11248.25
$ws.Range("L51").Value = 25000
$ws.Range("M51").Value = -10738.25
$ws.Range("N51").Value = -26020
$ws.Range("H61").Value = 47999.332
$ws.Range("I61").Value = 47999.332
$ws.Range("K61").Value = 47999.332
$ws.Range("M61").Value = -47707.332
$ws.Range("H75").Value = 0
$ws.Range("J75").Value = 0
$ws.Range("L75").Value = 0
$ws.Range("N75").ClearContents()
$ws.Range("H78").Value = 0
$ws.Range("J78").Value = 0
$ws.Range("L78").Value = 0
$ws.Range("N78").ClearContents()
$ws.Range("H100").Value = 869.5
$ws.Range("I100").Value = 869.5
$ws.Range("K100").Value = 1739
$ws.Range("M100").Value = -1198
$ws.Range("H107").Value = 1391
$ws.Range("I107").Value = 1391
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 4173
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -2253
$ws.Range("N107").ClearContents()
$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()
$ws.Range("H122").Value = 5376.0625
$ws.Range("I122").Value = 5785.923
$ws.Range("K122").Value = 17357.769
$ws.Range("M122").Value = -14907.769
$ws.Range("H136").Value = 324759.3
$ws.Range("I136").Value = 347046.2
$ws.Range("K136").Value = 1041138.6
$ws.Range("M136").Value = -1038588.6
